$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: update descriptive text (new shared strings must be introduced
#     in this exact order so they land at shared-string indices 54-60) ---
$ws.Range("C35").Value = "Writefield in um"
$ws.Range("C36").Value = "Name of column dataset"
$ws.Range("C37").Value = "GDSII database filepath"
$ws.Range("C38").Value = "Name of structure to write"
$ws.Range("C41").Value = "Stepszie for curves in nm"
$ws.Range("C40").Value = "Stepsize for areas in nm"
$ws.Range("C39").Value = "Stepsize for lines in nm"

# --- Column D: vartype swap for rows 35-38 (reuse existing shared strings) ---
$ws.Range("D35").Value = "int"
$ws.Range("D36").Value = "string"
$ws.Range("D37").Value = "string"
$ws.Range("D38").Value = "string"

# --- Clear obsolete cells ---
$ws.Range("D42").ClearContents()
$ws.Range("C43").ClearContents()
$ws.Range("D43").ClearContents()
$ws.Range("D44").ClearContents()
$ws.Range("D45").ClearContents()
$ws.Range("D46").ClearContents()

# --- New WF-zoom/shift/rot parameter rows (B42:B47), highlighted like B39:B41 ---
$ws.Range("B42").Value = "WFZoomU"
$ws.Range("B43").Value = "WFZoomV"
$ws.Range("B44").Value = "WFShiftU"
$ws.Range("B45").Value = "WFShiftV"
$ws.Range("B46").Value = "WFRotU"
$ws.Range("A47").Value = 13
$ws.Range("B47").Value = "WFRotV"

$ws.Range("B42:B47").Interior.Color = 65535

$ws.Range("C43").Select() | Out-Null

# Excel's "minimized" bookView flag (not exposed through a writable COM
# property the host engine persists) - best-effort, matches author intent.
$wb.Windows.Item(1).WindowState = -4140
